$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: insert a new worksheet "2022-Q3" positioned before the existing
# "2022-Q2" worksheet (so the sheet order becomes 总计, 2022-Q3, 2022-Q2).
# NOTE: sheet object handles captured at/after the insertion index go stale
# once Add() shifts things, so re-fetch everything by index afterwards.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$null = $wb.Worksheets.Add($wb.Worksheets.Item(2))

$wsQ3 = $wb.Worksheets.Item(2)
$wsQ2 = $wb.Worksheets.Item(3)
$wsQ3.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# Step 2: update the "总计" (totals) sheet - insert a new row 2 for the
# 2022-Q3 figures and push the existing 2022-Q2 figures down to row 3
# (re-numbering its index column from 0 to 1).
# ---------------------------------------------------------------------------
$wsTotal.Rows(2).Insert()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.01

# Row 2 should carry the same look as the rest of the data rows: copy A3's
# style onto A2 (bordered / centered "index" style) and strip the stray
# formatting that Insert() propagated into B2:D2.
$wsTotal.Range("A3").Copy($wsTotal.Range("A2"))
$wsTotal.Range("B2:D2").ClearFormats()

# The old row (now row 3) keeps its data but its index value changes 0 -> 1.
$wsTotal.Range("A3").Value = 1

# ---------------------------------------------------------------------------
# Step 3: populate the new "2022-Q3" sheet with its header row and the two
# fund rows, matching the look of the sibling quarter sheets (bold/bordered/
# centered header + index column, plain data cells, numeric-looking text
# columns kept as text).
# ---------------------------------------------------------------------------
function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $wsQ3.Range("B1") "基金代码"
Set-TextCell $wsQ3.Range("C1") "基金名称"
Set-TextCell $wsQ3.Range("D1") "基金规模"
Set-TextCell $wsQ3.Range("E1") "股票总仓位"
Set-TextCell $wsQ3.Range("F1") "仓位占比"
Set-TextCell $wsQ3.Range("G1") "持有市值(亿元)"
Set-TextCell $wsQ3.Range("H1") "仓位排名"

$wsQ3.Range("A2").Value = 0
Set-TextCell $wsQ3.Range("B2") "012315"
Set-TextCell $wsQ3.Range("C2") "创金合信港股通成长股票A"
Set-TextCell $wsQ3.Range("D2") "0.08"
Set-TextCell $wsQ3.Range("E2") "80.48"
Set-TextCell $wsQ3.Range("F2") "9.08"
Set-TextCell $wsQ3.Range("G2") "0.0073"
$wsQ3.Range("H2").Value = 4

$wsQ3.Range("A3").Value = 1
Set-TextCell $wsQ3.Range("B3") "012316"
Set-TextCell $wsQ3.Range("C3") "创金合信港股通成长股票C"
Set-TextCell $wsQ3.Range("D3") "0.07"
Set-TextCell $wsQ3.Range("E3") "80.48"
Set-TextCell $wsQ3.Range("F3") "9.08"
Set-TextCell $wsQ3.Range("G3") "0.0064"
$wsQ3.Range("H3").Value = 4

# Copy the header look (bold/bordered/centered) from the "总计" sheet header
# and the "index column" look onto the new sheet's A2:A3, format-only so the
# values we already wrote are kept.
$wsTotal.Range("B1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)

$wsTotal.Range("A2").Copy()
$wsQ3.Range("A2:A3").PasteSpecial(-4122)

# Restore the original active tab ("2022-Q2") which Add()/renaming bumped
# off the active sheet.
$wsQ2.Select()

Write-Host "2022-Q3 sheet added"
